$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update test_number column (A2:A11) from 4 to 3
$ws.Range("A2:A11").Value = 3

# Update the active selection to I12 (matches recorded selection in the diff)
$ws.Range("I12").Select()
